$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text "1" filled into several previously-empty inline-string cells
$ws.Range("U1").Value = "1"
$ws.Range("P9").Value = "1"
$ws.Range("H17").Value = "1"
$ws.Range("H19").Value = "1"
$ws.Range("T19").Value = "1"
$ws.Range("H21").Value = "1"
$ws.Range("T21").Value = "1"
$ws.Range("H23").Value = "1"
$ws.Range("H25").Value = "1"
$ws.Range("I69").Value = "1"

# Replace placeholder text "sn" / "123sad" with "1"
$ws.Range("A9").Value = "1"
$ws.Range("T17").Value = "1"
$ws.Range("A53").Value = "1"
$ws.Range("E59").Value = "1"

# Clear the "Fecha Recepción" date value in Q13, and drop its date number
# format back to the plain style used by its neighboring cells
$ws.Range("R13").Copy()
$ws.Range("Q13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Q13").Value = ""

# Mark "X" for STDF status checkboxes
$ws.Range("W29").Value = "X"
$ws.Range("T30").Value = "X"
$ws.Range("T31").Value = "X"
$ws.Range("T32").Value = "X"
$ws.Range("W33").Value = "X"
$ws.Range("W34").Value = "X"
$ws.Range("W35").Value = "X"
$ws.Range("W36").Value = "X"
$ws.Range("W37").Value = "X"
$ws.Range("T38").Value = "X"
$ws.Range("T39").Value = "X"
$ws.Range("T40").Value = "X"
$ws.Range("T41").Value = "X"
$ws.Range("W44").Value = "X"
$ws.Range("T45").Value = "X"

# Numeric "1" values
$ws.Range("P41").Value = 1
$ws.Range("M45").Value = 1
